$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 826; Excel shifts rows 826:940 down to 827:941
# and copies formatting (including the date style on column D) from the row above.
$ws.Rows(826).Insert()

# Populate the newly inserted row 826 with the new data point.
$ws.Cells.Item(826, 1).Value = 6
$ws.Cells.Item(826, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(826, 3).Value = "Metropolitana"
$ws.Cells.Item(826, 4).Value = 44918
$ws.Cells.Item(826, 5).Value = 13
$ws.Cells.Item(826, 6).Value = 100112003
$ws.Cells.Item(826, 7).Value = "Ajo"
$ws.Cells.Item(826, 8).Value = "Chino"
$ws.Cells.Item(826, 9).Value = "Primera"
$ws.Cells.Item(826, 10).Value = 1200
$ws.Cells.Item(826, 11).Value = 12000
$ws.Cells.Item(826, 12).Value = 12500
$ws.Cells.Item(826, 13).Value = 12208
$ws.Cells.Item(826, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(826, 15).Value = "China"
$ws.Cells.Item(826, 16).Value = 1221
$ws.Cells.Item(826, 17).Value = 10
$ws.Cells.Item(826, 18).Value = "Hortaliza"
